$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 26 (shifts old rows 26..52 down to 27..53)
$ws.Rows(26).Insert()

# Populate the new row 26 with the latest weekly price data point
$ws.Range("A26").Value = 3
$ws.Range("B26").Value = "Femacal de La Calera"
$ws.Range("C26").Value = "Coquimbo"
$ws.Range("D26").Value = 44894
$ws.Range("E26").Value = 5
$ws.Range("F26").Value = 300000000
$ws.Range("G26").Value = "Espárragos"
$ws.Range("H26").Value = "Verde"
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 1300
$ws.Range("K26").Value = 1500
$ws.Range("L26").Value = 1500
$ws.Range("M26").Value = 1500
$ws.Range("N26").Value = "$/kilo"
$ws.Range("O26").Value = "Provincia de Quillota"
$ws.Range("P26").Value = 1500
$ws.Range("Q26").Value = 1
$ws.Range("R26").Value = "Hortaliza"
